# members.xlsx update: refresh "Final" team roster (Sheet4) with each
# member's real LinkedIn profile link + photo filename (replacing the old
# placeholder github link / balidaan.jpg used for every row), re-wire the
# hyperlinks on column D to only the members that should carry a live link,
# and move the active sheet/selection to Sheet4.

$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item(4)

# --- Column D (link) / column E (img) per row ------------------------------
# row -> (linkedin url, image path)
$links = @{
    2  = "https://www.linkedin.com/in/vivek-gusain-656b78157/"
    3  = "https://www.linkedin.com/in/marmik07/"
    4  = "https://www.linkedin.com/in/dyknoww/"
    5  = "https://www.linkedin.com/in/mukulm03/"
    6  = "https://www.linkedin.com/in/sumit-sharma-a54814181/"
    7  = "https://www.linkedin.com/in/abhiraj-singh-rathore-170499/"
    8  = "https://www.linkedin.com/in/amol-bobade-97b34b179/"
    9  = "https://www.linkedin.com/in/divyanshu-bhaik-7438a6155/"
    10 = "https://www.linkedin.com/in/kunal-kishore-79742814b/"
    11 = "https://www.linkedin.com/in/sahaj-kulshrestha/"
    12 = "https://www.linkedin.com/in/parthivi-jain/"
    13 = "https://www.linkedin.com/in/varan-singh-rohila/"
    14 = "https://www.linkedin.com/in/achyut-sharma-7508a0b4/"
    15 = "https://www.linkedin.com/in/priyanka-kumar-a1135a1a1/"
    16 = "https://www.linkedin.com/in/rishi-kumar-5b808b153/"
}

$imgs = @{
    2  = "../members/vivek.jpg"
    3  = "../members/marmik.jpg"
    4  = "../members/mayank.jpg"
    5  = "../members/mukul.jpg"
    6  = "../members/sumit.jpg"
    7  = "../members/abhiraj.jpg"
    8  = "../members/amol.jpg"
    9  = "../members/divyanshu.png"
    10 = "../members/kunal.jpg"
    11 = "../members/sahaj.jpg"
    12 = "../members/parthivi.jpg"
    13 = "../members/varan.jpg"
    14 = "../members/achyut.jpg"
    15 = "../members/priyanka.jpg"
    16 = "../members/rishi.jpg"
}

for ($r = 2; $r -le 16; $r++) {
    $ws4.Cells.Item($r, 4).Value = $links[$r]
    $ws4.Cells.Item($r, 5).Value = $imgs[$r]
}

# --- Rebuild the hyperlinks on column D -------------------------------------
# Only these rows carry a live hyperlink object (the rest keep plain text).
$ws4.Hyperlinks.Delete()

$hyperlinkRows = @(4, 6, 7, 8, 9, 10, 13, 14, 15)
foreach ($r in $hyperlinkRows) {
    $addr = $links[$r]
    $ws4.Hyperlinks.Add($ws4.Cells.Item($r, 4), $addr, "", "", $addr)
}

# --- Active sheet / selection ----------------------------------------------
# Sheet4 ("Final" roster) becomes the active tab, selection moves to H33.
$ws4.Activate()
$ws4.Range("H33").Select()
